# Insert a new slide "Problem Solving Approach" at position 3 (right after
# "Objective and Scope" and before "Dataset Overview"), pushing all the
# subsequent slides down by one. Everything else in the deck is unchanged.

$p = $ppt.ActivePresentation

# "Title and Content" is CustomLayout #2 on this deck's slide master, and is
# the layout already used by the sibling content slides.
$s = $p.Slides.Add(3, 2)

# --- Title placeholder -----------------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Problem Solving Approach"

# --- Body / content placeholder ---------------------------------------
$body = $s.Shapes.Item(2)

# Position & size of the content placeholder (EMU -> points, 12700 EMU/pt).
$body.Left = 1443491 / 12700
$body.Top = 1759974 / 12700
$body.Width = 6571343 / 12700
$body.Height = 4293505 / 12700

$lines = @(
    "Here are the steps involved",
    "1) Load the dataset",
    "2) Perform Exploratory Data Analysis to extract valuable business insights",
    "3) Using NLTK for Data Pre-Processing and Feature Engineering.",
    "4) Transform Textual Descriptions to numerical features using techniques like TF-IDF or Word2Vec.",
    "5) Perform train-test split",
    "6) Train and evaluate various suitable Machine Learning and Deep Learning models.",
    "7) Get our predictions.",
    "8) Finally, we will compare the various models to get the best performing model."
)

$bodyText = [string]::Join("`r", $lines)

$tr = $body.TextFrame.TextRange
$tr.Text = $bodyText
$tr.Font.Size = 18

$body.TextFrame.AutoSize = 2
